# Daily attendance processing - 2025-10-09 07:42:20
# Reorders the "Recorded By" (column G) text for a specific set of rows:
#   "backup@backdoor.com, system, System" -> "system, backup@backdoor.com, System"
#   "dnasr281@gmail.com, System"          -> "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows whose "Recorded By" (column G) value is "backup@backdoor.com, system, System"
$rowsThreeParty = @(2, 29, 56)

foreach ($r in $rowsThreeParty) {
    $cell = $ws.Cells.Item($r, 7)
    $cell.Value = "system, backup@backdoor.com, System"
}

# Rows whose "Recorded By" (column G) value is "dnasr281@gmail.com, System"
$rowsTwoParty = @(3, 6, 11, 12, 13, 14, 15, 30, 33, 38, 39, 40, 41, 42, 57, 60, 65, 66, 67, 68, 69, 86, 89, 93, 95, 112, 115, 119, 121, 138, 141, 145, 147)

foreach ($r in $rowsTwoParty) {
    $cell = $ws.Cells.Item($r, 7)
    $cell.Value = "System, dnasr281@gmail.com"
}
